$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws_ALC.Range("H40").Value = 7008.5
$ws_ALC.Range("I40").Value = 10188.333
$ws_ALC.Range("K40").Value = 10188.333
$ws_ALC.Range("M40").Value = -10013.333

# ALC row 69
$ws_ALC.Range("H69").Value = 3805.2
$ws_ALC.Range("I69").Value = 3756.5
$ws_ALC.Range("J69").Value = 3837.6667
$ws_ALC.Range("K69").Value = 11269.5
$ws_ALC.Range("L69").Value = 11513.0001
$ws_ALC.Range("M69").Value = -10395.5
$ws_ALC.Range("N69").Value = -13261.0001

# ALC row 72
$ws_ALC.Range("H72").Value = 3805.2
$ws_ALC.Range("I72").Value = 3756.5
$ws_ALC.Range("J72").Value = 3837.6667
$ws_ALC.Range("K72").Value = 33808.5
$ws_ALC.Range("L72").Value = 34539.0003
$ws_ALC.Range("M72").Value = -29440.5
$ws_ALC.Range("N72").Value = -43275.0003

# ALC row 76
$ws_ALC.Range("H76").Value = 3267.1428
$ws_ALC.Range("I76").Value = 3003
$ws_ALC.Range("J76").Value = 3294.9473
$ws_ALC.Range("K76").Value = 3003
$ws_ALC.Range("L76").Value = 3294.9473
$ws_ALC.Range("M76").Value = -2688
$ws_ALC.Range("N76").Value = -3924.9473

# ALC row 79
$ws_ALC.Range("H79").Value = 3267.1428
$ws_ALC.Range("I79").Value = 3003
$ws_ALC.Range("J79").Value = 3294.9473
$ws_ALC.Range("K79").Value = 3003
$ws_ALC.Range("L79").Value = 3294.9473
$ws_ALC.Range("M79").Value = -1911
$ws_ALC.Range("N79").Value = -5478.9473

# ALC row 80
$ws_ALC.Range("H80").Value = 12821103
$ws_ALC.Range("I80").Value = 30303852
$ws_ALC.Range("J80").Value = 420.13333
$ws_ALC.Range("K80").Value = 90911556
$ws_ALC.Range("L80").Value = 1260.39999
$ws_ALC.Range("M80").Value = -90910558
$ws_ALC.Range("N80").Value = -3256.39999

# ALC row 83
$ws_ALC.Range("H83").Value = 12821103
$ws_ALC.Range("I83").Value = 30303852
$ws_ALC.Range("J83").Value = 420.13333
$ws_ALC.Range("K83").Value = 272734668
$ws_ALC.Range("L83").Value = 3781.19997
$ws_ALC.Range("M83").Value = -272729676
$ws_ALC.Range("N83").Value = -13765.19997

# ALC row 137
$ws_ALC.Range("H137").Value = 3749.1843
$ws_ALC.Range("I137").Value = 833.0833
$ws_ALC.Range("J137").Value = 8748.214
$ws_ALC.Range("K137").Value = 2499.2499
$ws_ALC.Range("L137").Value = 26244.642
$ws_ALC.Range("M137").Value = 50.7501000000002
$ws_ALC.Range("N137").Value = -31344.642

# ARM row 19
$ws_ARM.Range("H19").Value = 0
$ws_ARM.Range("J19").Value = 0
$ws_ARM.Range("L19").Value = 0
$ws_ARM.Range("N19").ClearContents()

# BSM row 107
$ws_BSM.Range("H107").Value = 1849.9642
$ws_BSM.Range("I107").Value = 1743.6471
$ws_BSM.Range("J107").Value = 2014.2727
$ws_BSM.Range("K107").Value = 1743.6471
$ws_BSM.Range("L107").Value = 2014.2727
$ws_BSM.Range("M107").Value = 176.3529000000001
$ws_BSM.Range("N107").Value = -5854.2727

# CUL row 131
$ws_CUL.Range("H131").Value = 877.87
$ws_CUL.Range("I131").Value = 500
$ws_CUL.Range("J131").Value = 881.6869
$ws_CUL.Range("K131").Value = 1500
$ws_CUL.Range("L131").Value = 2645.0607
$ws_CUL.Range("M131").Value = 3540
$ws_CUL.Range("N131").Value = -12725.0607

# GSM row 70
$ws_GSM.Range("H70").Value = 5213.391
$ws_GSM.Range("I70").Value = 5258.8423
$ws_GSM.Range("J70").Value = 4997.5
$ws_GSM.Range("K70").Value = 5258.8423
$ws_GSM.Range("L70").Value = 4997.5
$ws_GSM.Range("M70").Value = -4988.8423
$ws_GSM.Range("N70").Value = -5537.5

# GSM row 73
$ws_GSM.Range("H73").Value = 5213.391
$ws_GSM.Range("I73").Value = 5258.8423
$ws_GSM.Range("J73").Value = 4997.5
$ws_GSM.Range("K73").Value = 5258.8423
$ws_GSM.Range("L73").Value = 4997.5
$ws_GSM.Range("M73").Value = -4322.8423
$ws_GSM.Range("N73").Value = -6869.5

# GSM row 102
$ws_GSM.Range("H102").Value = 2003.7693
$ws_GSM.Range("I102").Value = 1898.619
$ws_GSM.Range("J102").Value = 2445.4
$ws_GSM.Range("K102").Value = 1898.619
$ws_GSM.Range("L102").Value = 2445.4
$ws_GSM.Range("M102").Value = -276.6189999999999
$ws_GSM.Range("N102").Value = -5689.4

# GSM row 122
$ws_GSM.Range("H122").Value = 1878.4348
$ws_GSM.Range("I122").Value = 1827.0667
$ws_GSM.Range("J122").Value = 1974.75
$ws_GSM.Range("K122").Value = 5481.2001
$ws_GSM.Range("L122").Value = 5924.25
$ws_GSM.Range("M122").Value = -3031.2001
$ws_GSM.Range("N122").Value = -10824.25

# LTW row 7
$ws_LTW.Range("H7").Value = 3492.9285
$ws_LTW.Range("I7").Value = 2520.4
$ws_LTW.Range("J7").Value = 4033.2222
$ws_LTW.Range("K7").Value = 2520.4
$ws_LTW.Range("L7").Value = 4033.2222
$ws_LTW.Range("M7").Value = -2408.4
$ws_LTW.Range("N7").Value = -4257.2222

# LTW row 40
$ws_LTW.Range("H40").Value = 5024.25
$ws_LTW.Range("I40").Value = 3675.4167
$ws_LTW.Range("J40").Value = 7047.5
$ws_LTW.Range("K40").Value = 3675.4167
$ws_LTW.Range("L40").Value = 7047.5
$ws_LTW.Range("M40").Value = -3539.4167
$ws_LTW.Range("N40").Value = -7319.5

# LTW row 61
$ws_LTW.Range("H61").Value = 4258.8
$ws_LTW.Range("I61").Value = 5164.6665
$ws_LTW.Range("J61").Value = 2900
$ws_LTW.Range("K61").Value = 5164.6665
$ws_LTW.Range("L61").Value = 2900
$ws_LTW.Range("M61").Value = -4962.6665
$ws_LTW.Range("N61").Value = -3304

# LTW row 82
$ws_LTW.Range("H82").Value = 5953652.5
$ws_LTW.Range("I82").Value = 1348.7778
$ws_LTW.Range("J82").Value = 16667799
$ws_LTW.Range("K82").Value = 1348.7778
$ws_LTW.Range("L82").Value = 16667799
$ws_LTW.Range("M82").Value = -987.7778000000001
$ws_LTW.Range("N82").Value = -16668521

# LTW row 85
$ws_LTW.Range("H85").Value = 5953652.5
$ws_LTW.Range("I85").Value = 1348.7778
$ws_LTW.Range("J85").Value = 16667799
$ws_LTW.Range("K85").Value = 1348.7778
$ws_LTW.Range("L85").Value = 16667799
$ws_LTW.Range("M85").Value = -100.7778000000001
$ws_LTW.Range("N85").Value = -16670295

# LTW row 93
$ws_LTW.Range("H93").Value = 1651.8096
$ws_LTW.Range("I93").Value = 1423
$ws_LTW.Range("J93").Value = 1792.6154
$ws_LTW.Range("K93").Value = 1423
$ws_LTW.Range("L93").Value = 1792.6154
$ws_LTW.Range("M93").Value = -175
$ws_LTW.Range("N93").Value = -4288.6154

# LTW row 113
$ws_LTW.Range("H113").Value = 4258.8
$ws_LTW.Range("I113").Value = 5164.6665
$ws_LTW.Range("J113").Value = 2900
$ws_LTW.Range("K113").Value = 5164.6665
$ws_LTW.Range("L113").Value = 2900
$ws_LTW.Range("M113").Value = -2994.6665
$ws_LTW.Range("N113").Value = -7240

# LTW row 126
$ws_LTW.Range("H126").Value = 3492.9285
$ws_LTW.Range("I126").Value = 2520.4
$ws_LTW.Range("J126").Value = 4033.2222
$ws_LTW.Range("K126").Value = 7561.200000000001
$ws_LTW.Range("L126").Value = 12099.6666
$ws_LTW.Range("M126").Value = -5091.200000000001
$ws_LTW.Range("N126").Value = -17039.6666

# WVR row 49
$ws_WVR.Range("H49").Value = 10000
$ws_WVR.Range("J49").Value = 10000
$ws_WVR.Range("L49").Value = 10000
$ws_WVR.Range("N49").Value = -10460

# WVR row 62
$ws_WVR.Range("H62").Value = 2999
$ws_WVR.Range("I62").Value = 0
$ws_WVR.Range("J62").Value = 2999
$ws_WVR.Range("K62").Value = 0
$ws_WVR.Range("L62").Value = 2999
$ws_WVR.Range("M62").ClearContents()
$ws_WVR.Range("N62").Value = -4247

# WVR row 65
$ws_WVR.Range("H65").Value = 2999
$ws_WVR.Range("I65").Value = 0
$ws_WVR.Range("J65").Value = 2999
$ws_WVR.Range("K65").Value = 0
$ws_WVR.Range("L65").Value = 14995
$ws_WVR.Range("M65").ClearContents()
$ws_WVR.Range("N65").Value = -21235

# WVR row 81
$ws_WVR.Range("H81").Value = 3137.4
$ws_WVR.Range("I81").Value = 850
$ws_WVR.Range("J81").Value = 7203.8887
$ws_WVR.Range("K81").Value = 1700
$ws_WVR.Range("L81").Value = 14407.7774
$ws_WVR.Range("M81").Value = -639
$ws_WVR.Range("N81").Value = -16529.7774

# WVR row 84
$ws_WVR.Range("H84").Value = 3137.4
$ws_WVR.Range("I84").Value = 850
$ws_WVR.Range("J84").Value = 7203.8887
$ws_WVR.Range("K84").Value = 8500
$ws_WVR.Range("L84").Value = 72038.887
$ws_WVR.Range("M84").Value = -3196
$ws_WVR.Range("N84").Value = -82646.887

# WVR row 122
$ws_WVR.Range("H122").Value = 3572229.5
$ws_WVR.Range("I122").Value = 4082405.2
$ws_WVR.Range("J122").Value = 1000
$ws_WVR.Range("K122").Value = 12247215.6
$ws_WVR.Range("L122").Value = 3000
$ws_WVR.Range("M122").Value = -12244765.6
$ws_WVR.Range("N122").Value = -7900

# WVR row 132
$ws_WVR.Range("H132").Value = 1649.1
$ws_WVR.Range("I132").Value = 1257.6
$ws_WVR.Range("J132").Value = 2823.6
$ws_WVR.Range("K132").Value = 3772.8
$ws_WVR.Range("L132").Value = 8470.799999999999
$ws_WVR.Range("M132").Value = -1242.8
$ws_WVR.Range("N132").Value = -13530.8
